{"js": "const replacements = [\n  [\"799\u00f73=266, 1\", \"326\u00f78=40, 6\"],\n  [\"579\u00f78=72, 3\", \"617\u00f73=205, 2\"],\n  [\"870\u00f78=108, 6\", \"333\u00f79=37, 0\"],\n  [\"455\u00f76=75, 5\", \"117\u00f77=16, 5\"],\n  [\"919\u00f77=131, 2\", \"400\u00f78=50, 0\"],\n  [\"425\u00f72=212, 1\", \"120\u00f73=40, 0\"],\n  [\"294\u00f78=36, 6\", \"987\u00f79=109, 6\"],\n  [\"741\u00f76=123, 3\", \"491\u00f79=54, 5\"],\n  [\"448\u00f77=64, 0\", \"114\u00f74=28, 2\"],\n  [\"284\u00f78=35, 4\", \"764\u00f78=95, 4\"],\n  [\"577\u00f76=96, 1\", \"337\u00f75=67, 2\"],\n  [\"135\u00f74=33, 3\", \"592\u00f78=74, 0\"],\n  [\"397\u00f73=132, 1\", \"123\u00f72=61, 1\"],\n  [\"467\u00f77=66, 5\", \"710\u00f75=142, 0\"],\n  [\"647\u00f77=92, 3\", \"681\u00f75=136, 1\"],\n  [\"375\u00f79=41, 6\", \"893\u00f73=297, 2\"],\n  [\"914\u00f74=228, 2\", \"413\u00f73=137, 2\"],\n  [\"585\u00f72=292, 1\", \"661\u00f74=165, 1\"],\n  [\"133\u00f75=26, 3\", \"989\u00f75=197, 4\"],\n  [\"207\u00f79=23, 0\", \"282\u00f72=141, 0\"],\n  [\"682\u00f77=97, 3\", \"808\u00f76=134, 4\"],\n  [\"321\u00f75=64, 1\", \"991\u00f76=165, 1\"],\n  [\"512\u00f72=256, 0\", \"775\u00f75=155, 0\"],\n  [\"703\u00f74=175, 3\", \"352\u00f76=58, 4\"],\n  [\"882\u00f76=147, 0\", \"310\u00f75=62, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$pairs = @(\n  @{ Old = \"799\u00f73=266, 1\"; New = \"326\u00f78=40, 6\" },\n  @{ Old = \"579\u00f78=72, 3\"; New = \"617\u00f73=205, 2\" },\n  @{ Old = \"870\u00f78=108, 6\"; New = \"333\u00f79=37, 0\" },\n  @{ Old = \"455\u00f76=75, 5\"; New = \"117\u00f77=16, 5\" },\n  @{ Old = \"919\u00f77=131, 2\"; New = \"400\u00f78=50, 0\" },\n  @{ Old = \"425\u00f72=212, 1\"; New = \"120\u00f73=40, 0\" },\n  @{ Old = \"294\u00f78=36, 6\"; New = \"987\u00f79=109, 6\" },\n  @{ Old = \"741\u00f76=123, 3\"; New = \"491\u00f79=54, 5\" },\n  @{ Old = \"448\u00f77=64, 0\"; New = \"114\u00f74=28, 2\" },\n  @{ Old = \"284\u00f78=35, 4\"; New = \"764\u00f78=95, 4\" },\n  @{ Old = \"577\u00f76=96, 1\"; New = \"337\u00f75=67, 2\" },\n  @{ Old = \"135\u00f74=33, 3\"; New = \"592\u00f78=74, 0\" },\n  @{ Old = \"397\u00f73=132, 1\"; New = \"123\u00f72=61, 1\" },\n  @{ Old = \"467\u00f77=66, 5\"; New = \"710\u00f75=142, 0\" },\n  @{ Old = \"647\u00f77=92, 3\"; New = \"681\u00f75=136, 1\" },\n  @{ Old = \"375\u00f79=41, 6\"; New = \"893\u00f73=297, 2\" },\n  @{ Old = \"914\u00f74=228, 2\"; New = \"413\u00f73=137, 2\" },\n  @{ Old = \"585\u00f72=292, 1\"; New = \"661\u00f74=165, 1\" },\n  @{ Old = \"133\u00f75=26, 3\"; New = \"989\u00f75=197, 4\" },\n  @{ Old = \"207\u00f79=23, 0\"; New = \"282\u00f72=141, 0\" },\n  @{ Old = \"682\u00f77=97, 3\"; New = \"808\u00f76=134, 4\" },\n  @{ Old = \"321\u00f75=64, 1\"; New = \"991\u00f76=165, 1\" },\n  @{ Old = \"512\u00f72=256, 0\"; New = \"775\u00f75=155, 0\" },\n  @{ Old = \"703\u00f74=175, 3\"; New = \"352\u00f76=58, 4\" },\n  @{ Old = \"882\u00f76=147, 0\"; New = \"310\u00f75=62, 0\" },\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $rng = $d.Content\n  $found = $rng.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n  if (-not $found) {\n    throw \"No match found for: $($pair.Old)\"\n  }\n}\n"}
